$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 102-105: "Bunny CDN" moves up to row 102, the rest shift down by one row
$ws.Range("B102").Value = "Bunny CDN"
$ws.Range("D102").Value = "Bunny_CDN.xml"

$ws.Range("B103").Value = "Cloudflare Blender"
$ws.Range("D103").Value = "Cloudflare_Blender_CDN.xml"

$ws.Range("B104").Value = "Google Cloud Storage"
$ws.Range("D104").Value = "Google Cloud Storage_CDN.xml"

$ws.Range("B105").Value = "Microsoft CDN"
$ws.Range("D105").Value = "Microsoft_CDN.xml"

# Rows 94-105: "In HTML" column flips from No to Yes
$ws.Range("E94:E105").Value = "Yes"

# Recalculate the workbook so formulas (H1/H2) and the chart reflect the new values
$excel.CalculateFullRebuild()
$wb.RefreshAll()

# Update the sheet view: the scroll position resets and the selection moves to the new last group
$ws.Range("E101:E105").Select()
